$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Add()
$ws.Name = "Table_1"
$ws.Range("A1").Value = "Hello"
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").Font.Italic = $true
$ws.Range("A1").Font.Size = 14
$ws.Columns.Item(1).ColumnWidth = 11.33
$ws.Rows.Item(1).RowHeight = 18
$ws.Range("A1:B1").Merge()
$ws.Range("A1").HorizontalAlignment = -4108
